$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.86"
$ws.Range("E2").Value = "'1.02%"
$ws.Range("D3").Value = "'27.26"
$ws.Range("E3").Value = "'0.92%"
$ws.Range("D4").Value = "'4.707"
$ws.Range("E4").Value = "'0.53%"
$ws.Range("D5").Value = "'0.06209"
$ws.Range("E5").Value = "'3.14%"
$ws.Range("E6").Value = "'0.65%"
$ws.Range("D7").Value = "'0.8505"
$ws.Range("E7").Value = "'-1.11%"
$ws.Range("D8").Value = "'0.9154"
$ws.Range("E8").Value = "'-1.27%"
$ws.Range("D9").Value = "'0.1411"
$ws.Range("E9").Value = "'0.86%"
$ws.Range("D10").Value = "'0.04647"
$ws.Range("E10").Value = "'-5.70%"
$ws.Range("D11").Value = "'0.07086"
$ws.Range("E11").Value = "'0.86%"
$ws.Range("D12").Value = "'0.03182"
$ws.Range("E12").Value = "'2.96%"
$ws.Range("D13").Value = "'0.09060"
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("D14").Value = "'0.001542"
$ws.Range("E14").Value = "'-0.22%"
$ws.Range("D15").Value = "'0.0006151"
$ws.Range("E15").Value = "'1.45%"
$ws.Range("D16").Value = "'0.006122"
$ws.Range("E16").Value = "'1.77%"
$ws.Range("D18").Value = "'3.170"
$ws.Range("E18").Value = "'-0.06%"
$ws.Range("E19").Value = "'0.64%"
$ws.Range("D21").Value = "'0.1311"
$ws.Range("E21").Value = "'0.97%"
$ws.Range("D22").Value = "'4.091"
$ws.Range("E22").Value = "'-0.98%"
$ws.Range("D23").Value = "'0.04243"
$ws.Range("E23").Value = "'0.43%"
$ws.Range("D24").Value = "'0.001213"
$ws.Range("E24").Value = "'-0.20%"
$ws.Range("D25").Value = "'0.004132"
$ws.Range("E25").Value = "'2.34%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.07%"
$ws.Range("E27").Value = "'5.05%"
$ws.Range("E40").Value = "'1.72%"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("E41").Value = "'-0.15%"
$ws.Range("D42").Value = "'0.004133"
$ws.Range("E42").Value = "'3.64%"
$ws.Range("D43").Value = "'0.002184"
$ws.Range("E43").Value = "'-0.76%"
$ws.Range("E44").Value = "'-7.44%"
$ws.Range("D45").Value = "'0.00005173"
$ws.Range("E45").Value = "'1.76%"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("E47").Value = "'-34.20%"
$ws.Range("D48").Value = "'0.1676"
$ws.Range("E48").Value = "'26.86%"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("E50").Value = "'0.04%"

Write-Host "Applied 61 cell updates"
